$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 214.77777
$ws.Range("I4").Value = 195.85715
$ws.Range("J4").Value = 281
$ws.Range("K4").Value = 195.85715
$ws.Range("L4").Value = 281
$ws.Range("M4").Value = -81.85714999999999
$ws.Range("N4").Value = -509
$ws.Range("H137").Value = 1636153.5
$ws.Range("I137").Value = 1840.2307
$ws.Range("J137").Value = 6947671.5
$ws.Range("K137").Value = 5520.6921
$ws.Range("L137").Value = 20843014.5
$ws.Range("M137").Value = -2970.6921
$ws.Range("N137").Value = -20848114.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5994.75
$ws.Range("I3").Value = 300
$ws.Range("J3").Value = 7893
$ws.Range("K3").Value = 300
$ws.Range("L3").Value = 7893
$ws.Range("M3").Value = -185
$ws.Range("N3").Value = -8123
$ws.Range("H6").Value = 223600.8
$ws.Range("J6").Value = 14000
$ws.Range("L6").Value = 14000
$ws.Range("N6").Value = -14346
$ws.Range("H60").Value = 14864.429
$ws.Range("I60").Value = 2025.5
$ws.Range("K60").Value = 2025.5
$ws.Range("M60").Value = -1292.5
$ws.Range("H62").Value = 29950
$ws.Range("J62").Value = 29950
$ws.Range("L62").Value = 29950
$ws.Range("N62").Value = -31198
$ws.Range("H65").Value = 29950
$ws.Range("J65").Value = 29950
$ws.Range("L65").Value = 89850
$ws.Range("N65").Value = -96090
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31622
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -98112
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H80").Value = 24000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H81").Value = 22666.666
$ws.Range("J81").Value = 22666.666
$ws.Range("L81").Value = 22666.666
$ws.Range("N81").Value = -24662.666
$ws.Range("H82").Value = 23500
$ws.Range("J82").Value = 23500
$ws.Range("L82").Value = 23500
$ws.Range("N82").Value = -24222
$ws.Range("H83").Value = 24000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H84").Value = 22666.666
$ws.Range("J84").Value = 22666.666
$ws.Range("L84").Value = 67999.99800000001
$ws.Range("N84").Value = -77983.99800000001
$ws.Range("H85").Value = 23500
$ws.Range("J85").Value = 23500
$ws.Range("L85").Value = 23500
$ws.Range("N85").Value = -25996
$ws.Range("H86").Value = 29700
$ws.Range("J86").Value = 29700
$ws.Range("L86").Value = 29700
$ws.Range("N86").Value = -32072
$ws.Range("H89").Value = 29700
$ws.Range("J89").Value = 29700
$ws.Range("L89").Value = 89100
$ws.Range("N89").Value = -100956
$ws.Range("H132").Value = 2996.9583
$ws.Range("I132").Value = 2154.923
$ws.Range("J132").Value = 3992.0908
$ws.Range("K132").Value = 6464.768999999999
$ws.Range("L132").Value = 11976.2724
$ws.Range("M132").Value = -3934.768999999999
$ws.Range("N132").Value = -17036.2724

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 25114.143
$ws.Range("I26").Value = 24000
$ws.Range("J26").Value = 25949.75
$ws.Range("K26").Value = 24000
$ws.Range("L26").Value = 25949.75
$ws.Range("M26").Value = -23708
$ws.Range("N26").Value = -26533.75
$ws.Range("H105").Value = 2341.9443
$ws.Range("I105").Value = 1928.8
$ws.Range("J105").Value = 2500.8462
$ws.Range("K105").Value = 1928.8
$ws.Range("L105").Value = 2500.8462
$ws.Range("M105").Value = -181.8
$ws.Range("N105").Value = -5994.8462

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2904.818
$ws.Range("I31").Value = 2154.7144
$ws.Range("K31").Value = 2154.7144
$ws.Range("M31").Value = -1859.7144
$ws.Range("H34").Value = 2904.818
$ws.Range("I34").Value = 2154.7144
$ws.Range("K34").Value = 2154.7144
$ws.Range("M34").Value = -1952.7144
$ws.Range("H35").Value = 359.2857
$ws.Range("I35").Value = 359.2857
$ws.Range("K35").Value = 359.2857
$ws.Range("M35").Value = -65.28570000000002
$ws.Range("H36").Value = 7663.2856
$ws.Range("I36").Value = 4728.8
$ws.Range("J36").Value = 14999.5
$ws.Range("K36").Value = 4728.8
$ws.Range("L36").Value = 14999.5
$ws.Range("M36").Value = -4340.8
$ws.Range("N36").Value = -15775.5
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H40").Value = 7663.2856
$ws.Range("I40").Value = 4728.8
$ws.Range("J40").Value = 14999.5
$ws.Range("K40").Value = 4728.8
$ws.Range("L40").Value = 14999.5
$ws.Range("M40").Value = -4568.8
$ws.Range("N40").Value = -15319.5
$ws.Range("H105").Value = 2355.5
$ws.Range("I105").Value = 2211
$ws.Range("K105").Value = 2211
$ws.Range("M105").Value = -464
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H118").Value = 32740
$ws.Range("J118").Value = 32740
$ws.Range("L118").Value = 32740
$ws.Range("N118").Value = -36054
$ws.Range("H123").Value = 33886.668
$ws.Range("J123").Value = 33886.668
$ws.Range("L123").Value = 33886.668
$ws.Range("N123").Value = -43686.668
$ws.Range("H132").Value = 2866.1082
$ws.Range("I132").Value = 2516.6924
$ws.Range("J132").Value = 3692
$ws.Range("K132").Value = 7550.0772
$ws.Range("L132").Value = 11076
$ws.Range("M132").Value = -5020.0772
$ws.Range("N132").Value = -16136

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1172.0526
$ws.Range("I4").Value = 95.57143000000001
$ws.Range("J4").Value = 1800
$ws.Range("K4").Value = 286.71429
$ws.Range("L4").Value = 5400
$ws.Range("M4").Value = -174.71429
$ws.Range("N4").Value = -5624
$ws.Range("H131").Value = 889.52
$ws.Range("I131").Value = 565
$ws.Range("J131").Value = 896.1429000000001
$ws.Range("K131").Value = 1695
$ws.Range("L131").Value = 2688.4287
$ws.Range("M131").Value = 3345
$ws.Range("N131").Value = -12768.4287

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 4789.8
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 4789.8
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 4789.8
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -5013.8
$ws.Range("H132").Value = 5524.1787
$ws.Range("I132").Value = 3843.4546
$ws.Range("J132").Value = 6611.706
$ws.Range("K132").Value = 11530.3638
$ws.Range("L132").Value = 19835.118
$ws.Range("M132").Value = -9000.363799999999
$ws.Range("N132").Value = -24895.118

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 55553.08
$ws.Range("I93").Value = 1791.9
$ws.Range("J93").Value = 234757
$ws.Range("K93").Value = 1791.9
$ws.Range("L93").Value = 234757
$ws.Range("M93").Value = -543.9000000000001
$ws.Range("N93").Value = -237253
$ws.Range("H100").Value = 4625
$ws.Range("I100").Value = 4700.4287
$ws.Range("J100").Value = 4566.3335
$ws.Range("K100").Value = 4700.4287
$ws.Range("L100").Value = 4566.3335
$ws.Range("M100").Value = -4159.4287
$ws.Range("N100").Value = -5648.3335
